# Add files via upload
# Updates the "Gen slack" sheet (adds a new row) and the "Lines" sheet
# (fixes the x_ohm_per_km value on row 3), then re-selects "Gen slack"
# as the active sheet/view.

$wb = $excel.ActiveWorkbook

# --- "Lines" sheet: correct F3 (x_ohm_per_km) from 0.083 to 0.0083 ---
$lines = $wb.Worksheets.Item("Lines")
$lines.Range("F3").Value = 0.0083
$lines.Range("J10").Select() | Out-Null

# --- "Gen slack" sheet: append a new row (A3 = 1) ---
$genSlack = $wb.Worksheets.Item("Gen slack")
$genSlack.Range("A3").Value = 1
$genSlack.Range("B3").Select() | Out-Null

# Activate "Gen slack" so it becomes the saved active tab
$genSlack.Activate() | Out-Null
